# Updates the cryptocurrency price/volume snapshot cells (columns D and E)
# with freshly scraped values, per the "Updated symbol list" GitHub Actions
# commit. Values are written with a leading apostrophe so Excel stores them
# as text (matching the workbook's original inline-string / text format for
# these numeric-looking price and percentage strings) and the cell style is
# then reset to "Normal" so no stray number-format/quote-prefix formatting
# is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.71%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'48.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.14%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.216"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.77%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07744"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.19%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.38%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.296"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'18.69%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-7.75%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-5.04%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1923"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.35%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09228"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.46%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04557"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.05%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.29%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001294"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.47%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'-1.58%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005889"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.43%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.344"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.62%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.402"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.34%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3440"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.12%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.087"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.39%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1367"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.30%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3036"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.001297"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.23%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004155"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.67%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.08%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003570"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-95.18%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02557"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-5.34%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05757"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'5.53%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01084"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'87.50%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007975"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.20%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1421"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.08%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008404"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'14.28%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007757"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-9.56%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3390"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'8.01%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006927"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.75%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000756"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'1.21%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05539"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-20.26%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004035"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'1.24%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002118"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'1.21%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002017"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'1.21%"
$ws.Range("E51").Style = "Normal"
